$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: the paragraph-mark run properties of the "I trust this email..."
# paragraph switch from <w:rFonts w:hint="eastAsia"/> to
# <w:rFonts w:hint="default"/>. That attribute isn't exposed through the
# regular Font/ParagraphFormat COM surface, so rebuild the paragraph via
# Range.InsertXML with the corrected pPr (the run content/formatting is kept
# identical to the source).
# ---------------------------------------------------------------------------
$marker = "I trust this email finds you well"

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Contains($marker)) {
        $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
               '<w:pPr><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr>' +
               '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr>' +
               '<w:t>I trust this email finds you well. As the coordinator for our project titled &quot;Matrix dissimilarity based on differences in moments and sparsity,&quot; I am reaching out to confirm the finalized author list.</w:t></w:r></w:p>'
        $p.Range.InsertXML($xml) | Out-Null
        break
    }
}

# ---------------------------------------------------------------------------
# Change 2: split the "Additionally, ..." run so that the new overleaf
# sentence is inserted, and move the "_GoBack" bookmark so that it sits
# between the two halves instead of sitting in its own paragraph further
# down (right after the "Thank you..." paragraph).
# ---------------------------------------------------------------------------
$oldSentence = "Additionally, should you have any revisions or comments, please communicate them before Berkeley time, 0:00 11th December."
$newSentence = "Additionally, should you have any revisions or comments, please communicate them or revise them directly in the overleaf, https://@@BOOKMARK_SPOT@@, before Berkeley time, 0:00 11th December."

$found = $d.Content.Find.Execute($oldSentence, $true, $false, $false, $false, $false, $true, 1, $false, $newSentence, 2)

# Locate the marker we just inserted, clear it, and drop the (moved) bookmark
# exactly at that now-empty spot.
$markerRange = $d.Content
$markerRange.Find.Execute("@@BOOKMARK_SPOT@@") | Out-Null
$markerRange.Text = ""
$d.Bookmarks.Add("_GoBack", $markerRange) | Out-Null
